$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TP#20088: Fixed smart post sheet
# Correct the tracking number in column E (rows 3-10) on the IMPB SmartPost sheet.
for ($r = 3; $r -le 10; $r++) {
    $ws.Range("E$r").Value = 630148367
}

# Rows 4-10 previously used a slightly different style (s="13") than row 3 (s="6").
# Copy row 3's formatting down so E4:E10 match E3's style, as in the fixed sheet.
$ws.Range("E3").Copy()
$ws.Range("E4:E10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reset the view: select E3:E10 (instead of the previous AO3 selection) and
# scroll back so the topLeftCell override is no longer needed.
$ws.Activate()
$ws.Range("E3:E10").Select() | Out-Null
